$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4" = -12.501
    "C6" = -12.318
    "C7" = -12.909
    "E7" = 16.385
    "C8" = -12.606
    "E11" = 17.324
    "E12" = 17.498
    "E15" = 16.188
    "C16" = -13.188
    "C20" = -12.205
    "E20" = 16.337
    "C21" = -12.725
    "E21" = 16.7
    "E22" = 16.646
    "E23" = 16.434
    "C28" = -12.848
    "C29" = -12.129
    "E29" = 16.916
    "C30" = -12.456
    "C32" = -12.785
    "E34" = 16.872
    "C40" = -12.445
    "E42" = 16.539
    "E43" = 16.977
    "E44" = 16.465
    "E45" = 16.667
    "C46" = -13.654
    "E46" = 16.956
    "E50" = 16.404
    "C51" = -11.257
    "E51" = 17.157
    "C52" = -11.355
    "C57" = -13.486
    "E57" = 16.603
    "C59" = -13.18
    "C62" = -13.555
    "E65" = 17.269
    "C66" = -11.579
    "E66" = 17.157
    "E67" = 17.314
    "C73" = -12.839
    "C74" = -11.866
    "C77" = -12.622
    "E79" = 16.862
    "E84" = 16.412
    "E87" = 16.501
    "C92" = -11.066
    "E92" = 17.769
    "E97" = 16.77
    "C100" = -12.726
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$wb.Save()